$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Columns.Item(1).ColumnWidth = 13.45

$ws.Range("A16").Value = "IPR Listo"
$ws.Range("A17").Value = "Proxy Listo"

$ws.Range("C11").Select()
Write-Output "done"
